$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2,3) {
    $ws.Range("D$row").Value = 0.0416
    $ws.Range("E$row").Value = 0.0883
    $ws.Range("K$row").Value = 16.2
    $ws.Range("L$row").Value = 0.3333333333333333
    $ws.Range("M$row").Value = 0
    $ws.Range("N$row").Value = 0
    $ws.Range("O$row").Value = 0
    $ws.Range("P$row").Value = 0
    $ws.Range("Q$row").Value = 0
    $ws.Range("R$row").Value = 0
    # buybacks_cash_returned (T) no longer populated for these rows
    $ws.Range("T$row").ClearContents()
    $ws.Range("U$row").Value = 29.6
    $ws.Range("V$row").Value = 0.1696275071633238
    $ws.Range("W$row").Value = 0.1062992125984252
    $ws.Range("X$row").Value = 0.04136624714725483
    $ws.Range("Y$row").Value = 0.06493296545117036
    $ws.Range("Z$row").Value = 0.3805496828752643
    $ws.Range("AB$row").Value = 0.04121983102901283
    $ws.Range("AC$row").Value = -0.04121983102901283
    $ws.Range("AD$row").Value = 1.49
    $ws.Range("AF$row").Value = 1.49
    $ws.Range("AG$row").Value = -28.11
    $ws.Range("AH$row").Value = 0.00846639013580317
    $ws.Range("AI$row").Value = 0.008151430603424693
    $ws.Range("AJ$row").Value = -0.1920213129312112
    $ws.Range("AK$row").Value = -0.1834976173379464
}
